# Update the NATMI LR-pairs table with newly recomputed TPM-based values.
# The underlying data was refreshed (new TPM normalization), which:
#  - changes most of the numeric columns (G..T) for the existing sending/target
#    cluster combinations,
#  - drops "MuSCs" as a possible Target cluster (column D) entirely, which
#    shrinks the table from 9 data rows (rows 2-10) down to 6 data rows
#    (rows 2-7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-7 (columns A-T). Column order matches the sheet header:
# A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E Ligand-expressing cells, F Ligand detection rate, G Ligand average
# expression value, H Ligand total expression value, I Ligand derived
# specificity of average expression value, J Ligand derived specificity of
# total expression value, K Receptor-expressing cells, L Receptor detection
# rate, M Receptor average expression value, N Receptor total expression
# value, O Receptor derived specificity of average expression value,
# P Receptor derived specificity of total expression value, Q Edge average
# expression weight, R Edge total expression weight, S Edge average
# expression derived specificity, T Edge total expression derived
# specificity.
$newData = @(
    @("ECs","Sema6d","Tyrobp","ECs",3,1,17.16209533333333,51.486286,0.2459970657298922,0.2459970657298922,2,0.6666666666666666,0.08160666666666666,0.24482,0.471786649605624,0.471786649605624,1.400541393168889,12.60487253852,0.1160581314535203,0.1160581314535203),
    @("ECs","Sema6d","Tyrobp","FAPs",3,1,17.16209533333333,51.486286,0.2459970657298922,0.2459970657298922,1,0.3333333333333333,0.09136699999999999,0.274101,0.528213350394376,0.528213350394376,1.568049164320666,14.112442478886,0.1299389342763719,0.1299389342763719),
    @("FAPs","Sema6d","Tyrobp","ECs",3,1,14.69090766666667,44.072723,0.2105756965403629,0.2105756965403628,2,0.6666666666666666,0.08160666666666666,0.24482,0.471786649605624,0.471786649605624,1.198876004984444,10.78988404486,0.09934680235914839,0.09934680235914836),
    @("FAPs","Sema6d","Tyrobp","FAPs",3,1,14.69090766666667,44.072723,0.2105756965403629,0.2105756965403628,1,0.3333333333333333,0.09136699999999999,0.274101,0.528213350394376,0.528213350394376,1.342264160780333,12.080377447023,0.1112288941812145,0.1112288941812145),
    @("MuSCs","Sema6d","Tyrobp","ECs",3,1,37.91244433333333,113.737333,0.543427237729745,0.543427237729745,2,0.6666666666666666,0.08160666666666666,0.24482,0.471786649605624,0.471786649605624,3.093908207228889,27.84517386506,0.2563817157929553,0.2563817157929553),
    @("MuSCs","Sema6d","Tyrobp","FAPs",3,1,37.91244433333333,113.737333,0.543427237729745,0.543427237729745,1,0.3333333333333333,0.09136699999999999,0.274101,0.528213350394376,0.528213350394376,3.463946301403666,31.175516712633,0.2870455219367897,0.2870455219367897)
)

for ($i = 0; $i -lt $newData.Length; $i++) {
    $rowIndex = $i + 2
    $rowValues = $newData[$i]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($rowIndex, $c + 1).Value = $rowValues[$c]
    }
}

# The refreshed data only has 6 rows (2-7); remove the now-obsolete rows
# 8-10 that corresponded to the dropped "MuSCs" target cluster.
$ws.Rows("8:10").Delete()
